# Update "想去人数" (column F) values on the "展览" sheet and the
# "全部类型" sheet, which both list the same events (the latter being
# the combined/aggregate view across all types).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    2  = 7162
    6  = 560
    7  = 176
    8  = 126
    16 = 1851
    18 = 5
    19 = 3729
    22 = 85
    23 = 35
    24 = 1
    25 = 32
    26 = 2376
    28 = 292
    31 = 40
    35 = 24
    37 = 1403
    38 = 136
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    2  = 7162
    7  = 560
    8  = 176
    9  = 126
    17 = 1851
    19 = 5
    20 = 3729
    23 = 85
    24 = 35
    25 = 1
    26 = 32
    27 = 2376
    29 = 292
    36 = 24
    38 = 1403
    39 = 136
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
